$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "62.896.38"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +2.92%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.976.79"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +1.61%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.01%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "594.86"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.38%  "

$ws.Range("E6").Value = "  +0.64%  "

$ws.Range("E7").Value = "  +0.02%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "2.975.89"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.55%  "

$ws.Range("E9").Value = "  +0.46%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "7.25"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +3.29%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.146"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +2.62%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.447"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.19%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000239"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +5.94%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "33.32"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.16%  "

$ws.Range("E15").Value = "  -0.36%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.465.84"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.54%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "62.793.31"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +2.76%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "6.73"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.01%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "2.972.85"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.50%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "443.41"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.23%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.49"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.12%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.673"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.03%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.10"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.19%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.31"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +2.09%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "81.86"
$ws.Range("D25").Style = "Normal"

$ws.Range("B26").Value = "InternetComputer(DFINITY)"
$ws.Range("C26").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "11.91"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.36%  "

$ws.Range("B27").Value = "Fetch.AI"
$ws.Range("C27").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.13"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -3.31%  "

$ws.Range("E28").Value = "  +0.09%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.22"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +3.56%  "

$ws.Range("E30").Value = "  -0.06%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.14"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -4.86%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.0₃0941"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +8.41%  "

$ws.Range("E33").Value = "  -0.08%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "26.76"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.12%  "

$ws.Range("E35").Value = "  +0.06%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.00"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.22%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.64"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.16%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.04"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.84%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.06"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +3.65%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "49.64"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.68%  "

$ws.Range("E41").Value = "  -0.34%  "

$ws.Range("E42").Value = "  -4.19%  "

$ws.Range("E43").Value = "  -0.64%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "40.06"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -4.71%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.747.17"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.31%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "135.43"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.32%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0341"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.54%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "364.48"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -2.24%  "

$ws.Range("B50").Value = "InjectiveProtocol"
$ws.Range("C50").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "23.09"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -3.14%  "

$ws.Range("B51").Value = "Stellar"
$ws.Range("C51").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.105"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.04%  "
